$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.748.70"
$ws.Range("E2").Value = "  -2.24%  "
$ws.Range("D3").Value = "1.875.19"
$ws.Range("E3").Value = "  -1.98%  "
$ws.Range("E4").Value = "  -0.89%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.17"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.688"
$ws.Range("E6").Value = "  -2.09%  "
$ws.Range("E7").Value = "  -0.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.40"
$ws.Range("E8").Value = "  +1.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.348"
$ws.Range("E9").Value = "  -2.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "50.50"
$ws.Range("E10").Value = "  -4.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0739"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("E12").Value = "  -2.21%  "
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("D14").Value = "2.147.02"
$ws.Range("E14").Value = "  -2.00%  "
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.89"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("D17").Value = "1.865.94"
$ws.Range("E17").Value = "  -2.77%  "
$ws.Range("D18").Value = "34.746.50"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.86"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "247.55"
$ws.Range("E21").Value = "  +2.04%  "
$ws.Range("E22").Value = "  -3.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.92"
$ws.Range("E23").Value = "  -3.24%  "
$ws.Range("E24").Value = "  -0.95%  "
$ws.Range("E25").Value = "  +3.37%  "
$ws.Range("E26").Value = "  -1.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.26"
$ws.Range("E27").Value = "  -1.86%  "
$ws.Range("E28").Value = "  -3.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.26"
$ws.Range("E29").Value = "  -2.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.128"
$ws.Range("E30").Value = "  -3.52%  "
$ws.Range("D31").Value = "4.128.61"
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("E32").Value = "  +14.47%  "
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  -1.69%  "
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("E37").Value = "  -5.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.834"
$ws.Range("E38").Value = "  -8.85%  "
$ws.Range("E39").Value = "  -3.31%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "98.16"
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.19"
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0658"
$ws.Range("E42").Value = "  +1.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0211"
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("E44").Value = "  -5.11%  "
$ws.Range("D45").Value = "1.292.93"
$ws.Range("E45").Value = "  -4.35%  "
$ws.Range("E46").Value = "  -4.19%  "
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("E48").Value = "  -2.17%  "
$ws.Range("E49").Value = "  +5.96%  "
$ws.Range("B50").Value = "Gas"
$ws.Range("C50").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "12.05"
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.49"
$ws.Range("E51").Value = "  -1.11%  "

Write-Host "Applied 84 cell updates"
